{"js": "// RELEASE: removed listed simulators from QRs and added reference to README.md\n//\n// 1) Remove the stray \"_GoBack\" bookmark that Word had left near the\n//    \"... functional parameters\" paragraph (an editing-position marker,\n//    not meaningful content).\n// 2) Replace the \"This VVC has been compiled and tested with Modelsim\n//    version 10.5b.\" paragraph with \"See README.md for a list of\n//    supported simulators.\" and drop its \"Liste\" list-paragraph style.\n// 3) Word re-drops a \"_GoBack\" bookmark at the new last-edit location,\n//    i.e. right after the replacement text.\n// 4) Bump the cached \"last update\" date in the footer from 2019-06-06 to\n//    2019-06-07.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- 1) Drop the old _GoBack bookmark (wherever Word last left it). ---\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) Find + replace the Modelsim sentence. ---\nconst target = body.search(\n  \"This VVC has been compiled and tested with Modelsim version 10.5b.\",\n  { matchCase: false }\n);\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  const hit = target.items[0];\n\n  // Swap the sentence.\n  hit.insertText(\"See README.md for a list of supported simulators.\", \"Replace\");\n\n  // The paragraph was styled \"Liste\" (list paragraph); the new sentence is\n  // plain body text, so reset the paragraph style back to the default.\n  const paras = hit.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n\n  const para = paras.items[0];\n  para.style = \"Normal\";\n\n  // --- 3) Re-plant _GoBack right after the new sentence. ---\n  const tail = para.getRange(\"End\");\n  tail.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 4) Update the cached DATE field result in the default footer. ---\nconst sections = doc.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const footer = section.getFooter(\"Primary\");\n  const dateHits = footer.search(\"2019-06-06\", { matchCase: true });\n  dateHits.load(\"items\");\n  await context.sync();\n\n  for (const hit of dateHits.items) {\n    hit.insertText(\"2019-06-07\", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"This VVC has been compiled and tested with Modelsim version 10.5b.\"\n$find.Execute() | Out-Null\n$range.Text = \"See README.md for a list of supported simulators.\"\n$range.Style = \"Normal\"\n"}
